$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move/resize the application window (mirrors the workbookView xWindow/yWindow
# change in the saved file: window was dragged further to the left).
$win = $wb.Windows.Item(1)
$win.Left = -28920
$win.Top = -120

# Update the cell value T2: 201884 -> 201634
$ws.Range("T2").Value = 201634

# Update the selected cell from T3 to T2
$ws.Range("T2").Select()
